$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 4930
$ws.Range("J10").Value = 4930
$ws.Range("L10").Value = 4930
$ws.Range("N10").Value = -5516

$ws.Range("H137").Value = 1221.1538
$ws.Range("I137").Value = 1073.7059
$ws.Range("J137").Value = 1499.6666
$ws.Range("K137").Value = 3221.1177
$ws.Range("L137").Value = 4498.9998
$ws.Range("M137").Value = -671.1176999999998
$ws.Range("N137").Value = -9598.9998

$ws.Range("H140").Value = 87564
$ws.Range("J140").Value = 87564
$ws.Range("L140").Value = 87564
$ws.Range("N140").Value = -97924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2318.2322
$ws.Range("I61").Value = 2292.9814
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2292.9814
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2080.9814
$ws.Range("N61").Value = -3424

$ws.Range("H97").Value = 2234.4583
$ws.Range("I97").Value = 1940.4706
$ws.Range("J97").Value = 2948.4285
$ws.Range("K97").Value = 1940.4706
$ws.Range("L97").Value = 2948.4285
$ws.Range("M97").Value = -1444.4706
$ws.Range("N97").Value = -3940.4285

$ws.Range("H110").Value = 1810.6471
$ws.Range("I110").Value = 1272.1818
$ws.Range("J110").Value = 2797.8333
$ws.Range("K110").Value = 1272.1818
$ws.Range("L110").Value = 2797.8333
$ws.Range("M110").Value = 772.8181999999999
$ws.Range("N110").Value = -6887.8333

$ws.Range("H136").Value = 2318.2322
$ws.Range("I136").Value = 2292.9814
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6878.9442
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4328.9442
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 734.6667
$ws.Range("I7").Value = 850
$ws.Range("J7").Value = 504
$ws.Range("K7").Value = 850
$ws.Range("L7").Value = 504
$ws.Range("M7").Value = -737
$ws.Range("N7").Value = -730

$ws.Range("H99").Value = 3060
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 4100
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 4100
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -7096

$ws.Range("H140").Value = 89700
$ws.Range("J140").Value = 89700
$ws.Range("L140").Value = 89700
$ws.Range("N140").Value = -100060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5521
$ws.Range("I31").Value = 4675.136
$ws.Range("J31").Value = 5998.154
$ws.Range("K31").Value = 4675.136
$ws.Range("L31").Value = 5998.154
$ws.Range("M31").Value = -4380.136
$ws.Range("N31").Value = -6588.154

$ws.Range("H34").Value = 5521
$ws.Range("I34").Value = 4675.136
$ws.Range("J34").Value = 5998.154
$ws.Range("K34").Value = 4675.136
$ws.Range("L34").Value = 5998.154
$ws.Range("M34").Value = -4473.136
$ws.Range("N34").Value = -6402.154

$ws.Range("H51").Value = 9350.667
$ws.Range("J51").Value = 9350.667
$ws.Range("L51").Value = 9350.667
$ws.Range("N51").Value = -10822.667

$ws.Range("H60").Value = 8250.5
$ws.Range("J60").Value = 8250.5
$ws.Range("L60").Value = 8250.5
$ws.Range("N60").Value = -9272.5

$ws.Range("H61").Value = 9350.667
$ws.Range("J61").Value = 9350.667
$ws.Range("L61").Value = 9350.667
$ws.Range("N61").Value = -10046.667

$ws.Range("H68").Value = 17200.428
$ws.Range("J68").Value = 17200.428
$ws.Range("L68").Value = 17200.428
$ws.Range("N68").Value = -18698.428

$ws.Range("H71").Value = 17200.428
$ws.Range("J71").Value = 17200.428
$ws.Range("L71").Value = 51601.284
$ws.Range("N71").Value = -59089.284

$ws.Range("H74").Value = 15224.5
$ws.Range("J74").Value = 17812.4
$ws.Range("L74").Value = 17812.4
$ws.Range("N74").Value = -19560.4

$ws.Range("H77").Value = 15224.5
$ws.Range("J77").Value = 17812.4
$ws.Range("L77").Value = 53437.2
$ws.Range("N77").Value = -62173.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1013.9545
$ws.Range("I113").Value = 850
$ws.Range("J113").Value = 1021.7619
$ws.Range("K113").Value = 2550
$ws.Range("L113").Value = 3065.2857
$ws.Range("M113").Value = -380
$ws.Range("N113").Value = -7405.2857

$ws.Range("H131").Value = 836.23
$ws.Range("I131").Value = 575
$ws.Range("J131").Value = 852.90424
$ws.Range("K131").Value = 1725
$ws.Range("L131").Value = 2558.71272
$ws.Range("M131").Value = 3315
$ws.Range("N131").Value = -12638.71272

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 7000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0

$ws.Range("H70").Value = 14430157
$ws.Range("I70").Value = 21640754
$ws.Range("J70").Value = 8962.923
$ws.Range("K70").Value = 21640754
$ws.Range("L70").Value = 8962.923
$ws.Range("M70").Value = -21640484
$ws.Range("N70").Value = -9502.923

$ws.Range("H73").Value = 14430157
$ws.Range("I73").Value = 21640754
$ws.Range("J73").Value = 8962.923
$ws.Range("K73").Value = 21640754
$ws.Range("L73").Value = 8962.923
$ws.Range("M73").Value = -21639818
$ws.Range("N73").Value = -10834.923

$ws.Range("H88").Value = 48000
$ws.Range("J88").Value = 48000
$ws.Range("L88").Value = 48000
$ws.Range("N88").Value = -48902

$ws.Range("H91").Value = 48000
$ws.Range("J91").Value = 48000
$ws.Range("L91").Value = 48000
$ws.Range("N91").Value = -51120

$ws.Range("H136").Value = 21190.54
$ws.Range("J136").Value = 21190.54
$ws.Range("L136").Value = 63571.62
$ws.Range("N136").Value = -68671.62

$ws.Range("H140").Value = 99893
$ws.Range("J140").Value = 99893
$ws.Range("L140").Value = 99893
$ws.Range("N140").Value = -110253

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 533.7963
$ws.Range("I22").Value = 440.28204
$ws.Range("J22").Value = 776.93335
$ws.Range("K22").Value = 440.28204
$ws.Range("L22").Value = 776.93335
$ws.Range("M22").Value = -145.28204
$ws.Range("N22").Value = -1366.93335

$ws.Range("H27").Value = 533.7963
$ws.Range("I27").Value = 440.28204
$ws.Range("J27").Value = 776.93335
$ws.Range("K27").Value = 440.28204
$ws.Range("L27").Value = 776.93335
$ws.Range("M27").Value = -333.28204
$ws.Range("N27").Value = -990.93335

$ws.Range("H46").Value = 900
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -712

$ws.Range("H132").Value = 3324.9092
$ws.Range("I132").Value = 4221.421
$ws.Range("J132").Value = 2108.2144
$ws.Range("K132").Value = 12664.263
$ws.Range("L132").Value = 6324.6432
$ws.Range("M132").Value = -10134.263
$ws.Range("N132").Value = -11384.6432

$ws.Range("H136").Value = 4336.5
$ws.Range("I136").Value = 3281.8333
$ws.Range("J136").Value = 5127.5
$ws.Range("K136").Value = 9845.499899999999
$ws.Range("L136").Value = 15382.5
$ws.Range("M136").Value = -7295.499899999999
$ws.Range("N136").Value = -20482.5

# Remove cells that should no longer have values
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M58").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N46").ClearContents()
